$d = $word.ActiveDocument

# Title line: date update
$d.Content.Find.Execute("2024-02-07 Wednesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-02-08 Thursday", 2)

$tbl = $d.Tables(1)

# Row 1
$tbl.Cell(1, 1).Range.Text = "77÷2="
$tbl.Cell(1, 2).Range.Text = "77÷2="
$tbl.Cell(1, 3).Range.Text = "12÷8="
$tbl.Cell(1, 4).Range.Text = "58÷9="
$tbl.Cell(1, 5).Range.Text = "87÷5="

# Row 5
$tbl.Cell(5, 1).Range.Text = "25÷3="
$tbl.Cell(5, 2).Range.Text = "23÷4="
$tbl.Cell(5, 3).Range.Text = "30÷5="
$tbl.Cell(5, 4).Range.Text = "95÷8="
$tbl.Cell(5, 5).Range.Text = "56÷2="

# Row 9
$tbl.Cell(9, 1).Range.Text = "25÷3="
$tbl.Cell(9, 2).Range.Text = "98÷4="
$tbl.Cell(9, 3).Range.Text = "97÷4="
$tbl.Cell(9, 4).Range.Text = "52÷6="
$tbl.Cell(9, 5).Range.Text = "75÷6="

# Row 13
$tbl.Cell(13, 1).Range.Text = "52÷8="
$tbl.Cell(13, 2).Range.Text = "83÷7="
$tbl.Cell(13, 3).Range.Text = "53÷4="
$tbl.Cell(13, 4).Range.Text = "63÷7="
$tbl.Cell(13, 5).Range.Text = "42÷8="

# Row 17
$tbl.Cell(17, 1).Range.Text = "96÷3="
$tbl.Cell(17, 2).Range.Text = "65÷3="
$tbl.Cell(17, 3).Range.Text = "66÷3="
$tbl.Cell(17, 4).Range.Text = "95÷8="
$tbl.Cell(17, 5).Range.Text = "77÷8="
